$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-33 down to 13-34.
$ws.Rows(12).Insert()

# Populate the new row 12 with the new record.
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C12").Value = "Los Lagos"
$ws.Range("D12").Value = 44614
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112030
$ws.Range("G12").Value = "Poroto granado"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 35000
$ws.Range("L12").Value = 35000
$ws.Range("M12").Value = 35000
$ws.Range("N12").Value = "$/saco 25 kilos"
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 1400
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
